# Add the "Horto Muso" team as a new row (row 35) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "Horto Muso"
$ws.Range("B35").Value = "Thomas Debiasi | MAI UNA GIOIA"
$ws.Range("C35").Value = "Giacomo Gasparini | MAI UNA GIOIA"
$ws.Range("D35").Value = "FEDERICO NICOLODI | U.S. Guarna"
$ws.Range("E35").Value = "Federico Andreis | iMontagna"
$ws.Range("F35").Value = "Daniele Dalbosco | SdrumALA"
